$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.306.46'
$ws.Range('E2').Value = '  +2.75%  '
$ws.Range('D3').Value = '1.903.03'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  -1.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.47'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5141'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3932'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08461'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.59'
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.269'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '1.901.96'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.71'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.346'
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.35'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06739'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.90'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.029'
$ws.Range('E22').Value = '  +0.95%  '
$ws.Range('D23').Value = '29.320.83'
$ws.Range('E23').Value = '  +2.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.17'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.221'
$ws.Range('E25').Value = '  -2.39%  '
$ws.Range('D26').Value = '2.119.91'
$ws.Range('E26').Value = '  +1.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.43'
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.95'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.449'
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.70'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.062'
$ws.Range('E31').Value = '  +1.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1049'
$ws.Range('E32').Value = '  -1.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.168'
$ws.Range('E33').Value = '  +6.25%  '
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02482'
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06581'
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.075'
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2202'
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.234'
$ws.Range('E39').Value = '  +3.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.143'
$ws.Range('E40').Value = '  +1.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6507'
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.233'
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6070'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.684'
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.057'
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.230'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.84'
$ws.Range('E51').Value = '  +0.93%  '
